# The deck currently applies the "Integral" design (green/teal palette) to
# the slide master (ppt/theme/theme1.xml) while the notes master carries a
# separate, plain "Office Theme" palette (ppt/theme/theme2.xml). The target
# revision swaps the two palettes: the main deck becomes the default blue
# "Office Theme" colors, and the notes master becomes the green "Integral"
# colors. Only the 12-color theme color scheme (and the theme/clrScheme
# display names, which PowerPoint derives from the applied palette) differ
# between the two theme parts - the font scheme and format scheme (fills,
# lines, effects) are already identical, so re-coloring the theme accounts
# for the whole change.

function Convert-HexToOleRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# New "Office Theme" 12-slot color scheme (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) that the main presentation theme should adopt.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation

# Re-color the presentation's (only) design / slide master theme so the
# whole deck switches from the old "Integral" palette to the standard
# "Office Theme" palette.
$slideTheme = $p.SlideMaster.Theme
$slideColorScheme = $slideTheme.ThemeColorScheme
for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $slideColorScheme.Item($i).RGB = Convert-HexToOleRgb $officeThemeColors[$i - 1]
}

# Keep the notes master in step with the same theme object (this deck's
# notes master theme tracks the presentation theme), so notes pages follow
# the recolor as well.
$notesTheme = $p.NotesMaster.Theme
$notesColorScheme = $notesTheme.ThemeColorScheme
for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $notesColorScheme.Item($i).RGB = Convert-HexToOleRgb $officeThemeColors[$i - 1]
}
